$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the days refreshed crypto data.
# For D-column values that Excel would otherwise auto-parse as numbers (losing the
# original text formatting, e.g. "1.00" -> 1), force text entry by temporarily
# switching the cell to a text NumberFormat, then restoring the original Style so no
# visible formatting/style actually changes.

$ws.Cells.Item(2, 4).Value = "56.705.69"
$ws.Cells.Item(2, 5).Value = "  +0.25%  "

$ws.Cells.Item(3, 4).Value = "2.404.48"
$ws.Cells.Item(3, 5).Value = "  -3.50%  "

$ws.Cells.Item(4, 5).Value = "  +0.29%  "

$cell = $ws.Cells.Item(5, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "487.36"
$cell.Style = $origStyle
$ws.Cells.Item(5, 5).Value = "  -1.33%  "

$cell = $ws.Cells.Item(6, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "152.86"
$cell.Style = $origStyle
$ws.Cells.Item(6, 5).Value = "  -0.31%  "

$ws.Cells.Item(7, 5).Value = "  +0.26%  "

$ws.Cells.Item(8, 5).Value = "  +18.28%  "

$ws.Cells.Item(9, 4).Value = "2.421.42"
$ws.Cells.Item(9, 5).Value = "  -3.25%  "

$cell = $ws.Cells.Item(10, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.97"
$cell.Style = $origStyle
$ws.Cells.Item(10, 5).Value = "  +3.36%  "

$cell = $ws.Cells.Item(11, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0995"
$cell.Style = $origStyle
$ws.Cells.Item(11, 5).Value = "  +0.63%  "

$cell = $ws.Cells.Item(12, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.334"
$cell.Style = $origStyle
$ws.Cells.Item(12, 5).Value = "  -0.32%  "

$ws.Cells.Item(13, 5).Value = "  +1.15%  "

$ws.Cells.Item(14, 4).Value = "2.848.86"
$ws.Cells.Item(14, 5).Value = "  -2.60%  "

$ws.Cells.Item(15, 4).Value = "57.085.42"
$ws.Cells.Item(15, 5).Value = "  +0.66%  "

$cell = $ws.Cells.Item(16, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "20.74"
$cell.Style = $origStyle
$ws.Cells.Item(16, 5).Value = "  -3.12%  "

$ws.Cells.Item(17, 5).Value = "  -2.63%  "

$ws.Cells.Item(18, 4).Value = "2.425.51"
$ws.Cells.Item(18, 5).Value = "  -3.06%  "

$cell = $ws.Cells.Item(19, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.74"
$cell.Style = $origStyle
$ws.Cells.Item(19, 5).Value = "  +4.13%  "

$cell = $ws.Cells.Item(20, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "324.31"
$cell.Style = $origStyle
$ws.Cells.Item(20, 5).Value = "  +0.75%  "

$cell = $ws.Cells.Item(21, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.98"
$cell.Style = $origStyle
$ws.Cells.Item(21, 5).Value = "  -3.32%  "

$cell = $ws.Cells.Item(22, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $origStyle
$ws.Cells.Item(22, 5).Value = "  +0.07%  "

$cell = $ws.Cells.Item(23, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.97"
$cell.Style = $origStyle
$ws.Cells.Item(23, 5).Value = "  +1.34%  "

$cell = $ws.Cells.Item(24, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "58.16"
$cell.Style = $origStyle
$ws.Cells.Item(24, 5).Value = "  -1.41%  "

$cell = $ws.Cells.Item(25, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.407"
$cell.Style = $origStyle
$ws.Cells.Item(25, 5).Value = "  -0.88%  "

$cell = $ws.Cells.Item(26, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Cells.Item(26, 5).Value = "  +0.02%  "

$cell = $ws.Cells.Item(27, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.160"
$cell.Style = $origStyle
$ws.Cells.Item(27, 5).Value = "  -1.71%  "

$ws.Cells.Item(28, 4).Value = "2.537.20"
$ws.Cells.Item(28, 5).Value = "  -2.98%  "

$cell = $ws.Cells.Item(29, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.30"
$cell.Style = $origStyle
$ws.Cells.Item(29, 5).Value = "  -4.05%  "

$ws.Cells.Item(30, 4).Value = "0.0₃0782"
$ws.Cells.Item(30, 5).Value = "  -3.84%  "

$ws.Cells.Item(31, 5).Value = "  +0.24%  "

$cell = $ws.Cells.Item(32, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "149.96"
$cell.Style = $origStyle
$ws.Cells.Item(32, 5).Value = "  -1.18%  "

$cell = $ws.Cells.Item(33, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "18.55"
$cell.Style = $origStyle

$ws.Cells.Item(34, 5).Value = "  -0.05%  "

$cell = $ws.Cells.Item(35, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.34"
$cell.Style = $origStyle
$ws.Cells.Item(35, 5).Value = "  +1.71%  "

$cell = $ws.Cells.Item(36, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.15"
$cell.Style = $origStyle
$ws.Cells.Item(36, 5).Value = "  -1.15%  "

$ws.Cells.Item(37, 5).Value = "  -2.11%  "

$cell = $ws.Cells.Item(38, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.845"
$cell.Style = $origStyle
$ws.Cells.Item(38, 5).Value = "  -2.72%  "

$cell = $ws.Cells.Item(39, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.101"
$cell.Style = $origStyle
$ws.Cells.Item(39, 5).Value = "  +9.42%  "

$cell = $ws.Cells.Item(40, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "34.10"
$cell.Style = $origStyle
$ws.Cells.Item(40, 5).Value = "  +0.45%  "

$ws.Cells.Item(41, 5).Value = "  +0.59%  "

$cell = $ws.Cells.Item(42, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.36"
$cell.Style = $origStyle
$ws.Cells.Item(42, 5).Value = "  -2.44%  "

$cell = $ws.Cells.Item(43, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.995"
$cell.Style = $origStyle
$ws.Cells.Item(43, 5).Value = "  -0.01%  "

$cell = $ws.Cells.Item(44, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.592"
$cell.Style = $origStyle
$ws.Cells.Item(44, 5).Value = "  -4.12%  "

$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "268.93"
$cell.Style = $origStyle
$ws.Cells.Item(45, 5).Value = "  -0.12%  "

$cell = $ws.Cells.Item(46, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0532"
$cell.Style = $origStyle
$ws.Cells.Item(46, 5).Value = "  -5.81%  "

$cell = $ws.Cells.Item(47, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.21"
$cell.Style = $origStyle
$ws.Cells.Item(47, 5).Value = "  -0.05%  "

$ws.Cells.Item(48, 5).Value = "  -1.20%  "

$cell = $ws.Cells.Item(49, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.57"
$cell.Style = $origStyle
$ws.Cells.Item(49, 5).Value = "  -7.15%  "

$cell = $ws.Cells.Item(50, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "17.38"
$cell.Style = $origStyle
$ws.Cells.Item(50, 5).Value = "  -2.55%  "

$ws.Cells.Item(51, 4).Value = "1.865.72"
$ws.Cells.Item(51, 5).Value = "  -1.32%  "
